# Extend the age-lookup tables on the "HESD-FoHERbA" sheet.
#
# The sheet already has a header row (row 1) running B1:BT1 = 0..70 and a
# data row (row 2) running B2:BT2 = Calcs!C89..Calcs!BU89 (a "retired
# share" lookup curve that ends at 1 once the equipment is fully retired).
# The table only went up to age 70; extend both rows out to age 210
# (column HD) so lookups for older equipment resolve instead of
# #N/A-ing past the end of the table. The newly appended row-2 cells are
# plain "1" constants (the survival/retirement curve is already flat at
# 1 by that point), not formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HESD-FoHERbA")

$firstNewCol = 73   # column BU (existing data stops at BT = 72)
$lastNewCol  = 212  # column HD -> age 210

for ($col = $firstNewCol; $col -le $lastNewCol; $col++) {
    $ws.Cells.Item(1, $col).Value = $col - 2
    $ws.Cells.Item(2, $col).Value = 1
}

# Mirror the author's final view state: the "HESD-FoHERbA" tab ends up
# active/selected, scrolled over to the newly extended columns, with
# BS17 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 61
$ws.Range("BS17").Select()
